$d = $word.ActiveDocument

# The "Requisitos" list paragraph holds three requirement lines, each its own
# run ending in a manual line break (w:br, represented in Range.Text as
# Chr(11)). The edit moves the first line ("LOQ4083 - Fenomenos de
# Transporte I (Requisito fraco)") from the top of the list to the bottom,
# leaving the other two lines ("LOB1006 ...", "LOB1019 ...") untouched and
# in the same relative order.

$target = "LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)"

# Locate the run containing the text we need to move.
$moveRange = $d.Content.Duplicate
$found = $moveRange.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the text to move: $target"
}

# $moveRange now spans just the "LOQ4083 ..." text (Find.Execute with no
# replacement collapses the match onto moveRange). Extend it one character
# so it also swallows the following manual line break (Chr(11)).
$lineStart = $moveRange.Start
$lineBreakEnd = $moveRange.End + 1
$fullLineRange = $d.Range($lineStart, $lineBreakEnd)
$savedText = $fullLineRange.Text

# Find the paragraph that contains this line (walk the Paragraphs collection
# rather than Range.Paragraphs, which only "sees" the narrow sub-range).
$containingParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($lineStart -ge $candidate.Range.Start -and $lineStart -lt $candidate.Range.End) {
        $containingParagraph = $candidate
        break
    }
}

if ($null -eq $containingParagraph) {
    throw "Could not find the paragraph containing: $target"
}

# Remove the line (text + its trailing break) from its current position ...
$fullLineRange.Delete() | Out-Null

# ... and re-insert it just before the paragraph mark, i.e. after the last
# remaining line in the list.
$paraEnd = $containingParagraph.Range.End
$insertionPoint = $d.Range($paraEnd - 1, $paraEnd - 1)
$insertionPoint.InsertAfter($savedText) | Out-Null

Write-Output "Moved requirement line. Paragraph now reads: $($containingParagraph.Range.Text)"
